# Weekly update: a new price record is inserted as the second data row
# (row 8), shifting every following record down by one row and adding a
# single new record at the end (old last row becomes the new last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8 - this shifts rows 8:127 down to 9:128
# and extends the sheet dimension to A1:R128, exactly like the diff shows.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record's data.
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 45168
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112040
$ws.Range("G8").Value = "Cilantro"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1089
$ws.Range("N8").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 544
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = "Hortaliza"
